$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = "ALC"; Cell = "H46"; Value = 2555 }
    @{ Sheet = "ALC"; Cell = "J46"; Value = 2555 }
    @{ Sheet = "ALC"; Cell = "L46"; Value = 7665 }
    @{ Sheet = "ALC"; Cell = "N46"; Value = -7903 }
    @{ Sheet = "ALC"; Cell = "H58"; Value = 193.33333 }
    @{ Sheet = "ALC"; Cell = "J58"; Value = 180 }
    @{ Sheet = "ALC"; Cell = "L58"; Value = 540 }
    @{ Sheet = "ALC"; Cell = "N58"; Value = -840 }
    @{ Sheet = "ALC"; Cell = "H60"; Value = 2555 }
    @{ Sheet = "ALC"; Cell = "J60"; Value = 2555 }
    @{ Sheet = "ALC"; Cell = "L60"; Value = 7665 }
    @{ Sheet = "ALC"; Cell = "N60"; Value = -8633 }
    @{ Sheet = "ALC"; Cell = "H62"; Value = 0 }
    @{ Sheet = "ALC"; Cell = "I62"; Value = 0 }
    @{ Sheet = "ALC"; Cell = "K62"; Value = 0 }
    @{ Sheet = "ALC"; Cell = "M62"; Value = $null }
    @{ Sheet = "ALC"; Cell = "H65"; Value = 0 }
    @{ Sheet = "ALC"; Cell = "I65"; Value = 0 }
    @{ Sheet = "ALC"; Cell = "K65"; Value = 0 }
    @{ Sheet = "ALC"; Cell = "M65"; Value = $null }
    @{ Sheet = "ALC"; Cell = "H86"; Value = 224197.89 }
    @{ Sheet = "ALC"; Cell = "I86"; Value = 2175.4 }
    @{ Sheet = "ALC"; Cell = "K86"; Value = 2175.4 }
    @{ Sheet = "ALC"; Cell = "M86"; Value = -1052.4 }
    @{ Sheet = "ALC"; Cell = "H89"; Value = 224197.89 }
    @{ Sheet = "ALC"; Cell = "I89"; Value = 2175.4 }
    @{ Sheet = "ALC"; Cell = "K89"; Value = 10877 }
    @{ Sheet = "ALC"; Cell = "M89"; Value = -5261 }
    @{ Sheet = "ARM"; Cell = "H2"; Value = 632.3 }
    @{ Sheet = "ARM"; Cell = "I2"; Value = 526.5 }
    @{ Sheet = "ARM"; Cell = "K2"; Value = 526.5 }
    @{ Sheet = "ARM"; Cell = "M2"; Value = -413.5 }
    @{ Sheet = "ARM"; Cell = "H102"; Value = 4999 }
    @{ Sheet = "ARM"; Cell = "I102"; Value = 4999 }
    @{ Sheet = "ARM"; Cell = "J102"; Value = 0 }
    @{ Sheet = "ARM"; Cell = "K102"; Value = 4999 }
    @{ Sheet = "ARM"; Cell = "L102"; Value = 0 }
    @{ Sheet = "ARM"; Cell = "M102"; Value = $null }
    @{ Sheet = "ARM"; Cell = "N102"; Value = -3377 }
    @{ Sheet = "ARM"; Cell = "H116"; Value = 632.3 }
    @{ Sheet = "ARM"; Cell = "I116"; Value = 526.5 }
    @{ Sheet = "ARM"; Cell = "K116"; Value = 526.5 }
    @{ Sheet = "ARM"; Cell = "M116"; Value = 1767.5 }
    @{ Sheet = "ARM"; Cell = "H132"; Value = 2227.875 }
    @{ Sheet = "ARM"; Cell = "I132"; Value = 2635 }
    @{ Sheet = "ARM"; Cell = "J132"; Value = 1006.5 }
    @{ Sheet = "ARM"; Cell = "K132"; Value = 7905 }
    @{ Sheet = "ARM"; Cell = "L132"; Value = 3019.5 }
    @{ Sheet = "ARM"; Cell = "M132"; Value = -5375 }
    @{ Sheet = "ARM"; Cell = "N132"; Value = -8079.5 }
    @{ Sheet = "BSM"; Cell = "H3"; Value = 632.3 }
    @{ Sheet = "BSM"; Cell = "I3"; Value = 526.5 }
    @{ Sheet = "BSM"; Cell = "K3"; Value = 526.5 }
    @{ Sheet = "BSM"; Cell = "M3"; Value = -412.5 }
    @{ Sheet = "BSM"; Cell = "H94"; Value = 478 }
    @{ Sheet = "BSM"; Cell = "I94"; Value = 401.22223 }
    @{ Sheet = "BSM"; Cell = "K94"; Value = 401.22223 }
    @{ Sheet = "BSM"; Cell = "M94"; Value = 49.77776999999998 }
    @{ Sheet = "BSM"; Cell = "H99"; Value = 1808.25 }
    @{ Sheet = "BSM"; Cell = "I99"; Value = 1808.25 }
    @{ Sheet = "BSM"; Cell = "J99"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "K99"; Value = 1808.25 }
    @{ Sheet = "BSM"; Cell = "L99"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "M99"; Value = $null }
    @{ Sheet = "BSM"; Cell = "N99"; Value = -310.25 }
    @{ Sheet = "BSM"; Cell = "H105"; Value = 2041 }
    @{ Sheet = "BSM"; Cell = "I105"; Value = 1457.4 }
    @{ Sheet = "BSM"; Cell = "K105"; Value = 1457.4 }
    @{ Sheet = "BSM"; Cell = "M105"; Value = 289.5999999999999 }
    @{ Sheet = "BSM"; Cell = "H134"; Value = 2064 }
    @{ Sheet = "BSM"; Cell = "I134"; Value = 1990.9231 }
    @{ Sheet = "BSM"; Cell = "K134"; Value = 5972.7693 }
    @{ Sheet = "BSM"; Cell = "M134"; Value = -3437.7693 }
    @{ Sheet = "CRP"; Cell = "H4"; Value = 299 }
    @{ Sheet = "CRP"; Cell = "I4"; Value = 224 }
    @{ Sheet = "CRP"; Cell = "K4"; Value = 224 }
    @{ Sheet = "CRP"; Cell = "M4"; Value = -112 }
    @{ Sheet = "CRP"; Cell = "H7"; Value = 34.125 }
    @{ Sheet = "CRP"; Cell = "I7"; Value = 30.166666 }
    @{ Sheet = "CRP"; Cell = "J7"; Value = 36.5 }
    @{ Sheet = "CRP"; Cell = "K7"; Value = 30.166666 }
    @{ Sheet = "CRP"; Cell = "L7"; Value = 36.5 }
    @{ Sheet = "CRP"; Cell = "M7"; Value = 82.83333400000001 }
    @{ Sheet = "CRP"; Cell = "N7"; Value = -262.5 }
    @{ Sheet = "CRP"; Cell = "H31"; Value = 2282.7778 }
    @{ Sheet = "CRP"; Cell = "I31"; Value = 2372 }
    @{ Sheet = "CRP"; Cell = "J31"; Value = 2104.3333 }
    @{ Sheet = "CRP"; Cell = "K31"; Value = 2372 }
    @{ Sheet = "CRP"; Cell = "L31"; Value = 2104.3333 }
    @{ Sheet = "CRP"; Cell = "M31"; Value = -2077 }
    @{ Sheet = "CRP"; Cell = "N31"; Value = -2694.3333 }
    @{ Sheet = "CRP"; Cell = "H32"; Value = 1211 }
    @{ Sheet = "CRP"; Cell = "I32"; Value = 0 }
    @{ Sheet = "CRP"; Cell = "J32"; Value = 1211 }
    @{ Sheet = "CRP"; Cell = "K32"; Value = 0 }
    @{ Sheet = "CRP"; Cell = "L32"; Value = $null }
    @{ Sheet = "CRP"; Cell = "M32"; Value = 1211 }
    @{ Sheet = "CRP"; Cell = "N32"; Value = -1843 }
    @{ Sheet = "CRP"; Cell = "H34"; Value = 2282.7778 }
    @{ Sheet = "CRP"; Cell = "I34"; Value = 2372 }
    @{ Sheet = "CRP"; Cell = "J34"; Value = 2104.3333 }
    @{ Sheet = "CRP"; Cell = "K34"; Value = 2372 }
    @{ Sheet = "CRP"; Cell = "L34"; Value = 2104.3333 }
    @{ Sheet = "CRP"; Cell = "M34"; Value = -2170 }
    @{ Sheet = "CRP"; Cell = "N34"; Value = -2508.3333 }
    @{ Sheet = "CRP"; Cell = "H37"; Value = 9999 }
    @{ Sheet = "CRP"; Cell = "J37"; Value = 9999 }
    @{ Sheet = "CRP"; Cell = "L37"; Value = 9999 }
    @{ Sheet = "CRP"; Cell = "N37"; Value = -10213 }
    @{ Sheet = "CRP"; Cell = "H41"; Value = 1000 }
    @{ Sheet = "CRP"; Cell = "I41"; Value = 1000 }
    @{ Sheet = "CRP"; Cell = "K41"; Value = 1000 }
    @{ Sheet = "CRP"; Cell = "M41"; Value = -572 }
    @{ Sheet = "CRP"; Cell = "H58"; Value = 2169.4 }
    @{ Sheet = "CRP"; Cell = "I58"; Value = 1833.3334 }
    @{ Sheet = "CRP"; Cell = "J58"; Value = 2673.5 }
    @{ Sheet = "CRP"; Cell = "K58"; Value = 1833.3334 }
    @{ Sheet = "CRP"; Cell = "L58"; Value = 2673.5 }
    @{ Sheet = "CRP"; Cell = "M58"; Value = -1630.3334 }
    @{ Sheet = "CRP"; Cell = "N58"; Value = -3079.5 }
    @{ Sheet = "CRP"; Cell = "H60"; Value = 19000 }
    @{ Sheet = "CRP"; Cell = "J60"; Value = 19000 }
    @{ Sheet = "CRP"; Cell = "L60"; Value = 19000 }
    @{ Sheet = "CRP"; Cell = "N60"; Value = -20022 }
    @{ Sheet = "CRP"; Cell = "H122"; Value = 2509.6 }
    @{ Sheet = "CRP"; Cell = "I122"; Value = 1937 }
    @{ Sheet = "CRP"; Cell = "K122"; Value = 5811 }
    @{ Sheet = "CRP"; Cell = "M122"; Value = -3361 }
    @{ Sheet = "CRP"; Cell = "H136"; Value = 2169.4 }
    @{ Sheet = "CRP"; Cell = "I136"; Value = 1833.3334 }
    @{ Sheet = "CRP"; Cell = "J136"; Value = 2673.5 }
    @{ Sheet = "CRP"; Cell = "K136"; Value = 5500.0002 }
    @{ Sheet = "CRP"; Cell = "L136"; Value = 8020.5 }
    @{ Sheet = "CRP"; Cell = "M136"; Value = -2950.0002 }
    @{ Sheet = "CRP"; Cell = "N136"; Value = -13120.5 }
    @{ Sheet = "CUL"; Cell = "H2"; Value = 87.34999999999999 }
    @{ Sheet = "CUL"; Cell = "J2"; Value = 138.36363 }
    @{ Sheet = "CUL"; Cell = "L2"; Value = 830.18178 }
    @{ Sheet = "CUL"; Cell = "N2"; Value = -1056.18178 }
    @{ Sheet = "CUL"; Cell = "H45"; Value = 1416.5 }
    @{ Sheet = "CUL"; Cell = "I45"; Value = 800 }
    @{ Sheet = "CUL"; Cell = "J45"; Value = 2033 }
    @{ Sheet = "CUL"; Cell = "K45"; Value = 2400 }
    @{ Sheet = "CUL"; Cell = "L45"; Value = 6099 }
    @{ Sheet = "CUL"; Cell = "M45"; Value = -1868 }
    @{ Sheet = "CUL"; Cell = "N45"; Value = -7163 }
    @{ Sheet = "CUL"; Cell = "H94"; Value = 15497.143 }
    @{ Sheet = "CUL"; Cell = "J94"; Value = 15497.143 }
    @{ Sheet = "CUL"; Cell = "L94"; Value = 46491.429 }
    @{ Sheet = "CUL"; Cell = "N94"; Value = -47843.429 }
    @{ Sheet = "CUL"; Cell = "H116"; Value = 489.5 }
    @{ Sheet = "CUL"; Cell = "I116"; Value = 489.5 }
    @{ Sheet = "CUL"; Cell = "K116"; Value = 1468.5 }
    @{ Sheet = "CUL"; Cell = "M116"; Value = 1973.5 }
    @{ Sheet = "CUL"; Cell = "H131"; Value = 2512.3333 }
    @{ Sheet = "CUL"; Cell = "I131"; Value = 1537.8 }
    @{ Sheet = "CUL"; Cell = "J131"; Value = 2999.6 }
    @{ Sheet = "CUL"; Cell = "K131"; Value = 4613.4 }
    @{ Sheet = "CUL"; Cell = "L131"; Value = 8998.799999999999 }
    @{ Sheet = "CUL"; Cell = "M131"; Value = 426.6000000000004 }
    @{ Sheet = "CUL"; Cell = "N131"; Value = -19078.8 }
    @{ Sheet = "GSM"; Cell = "H70"; Value = 2801 }
    @{ Sheet = "GSM"; Cell = "I70"; Value = 2801 }
    @{ Sheet = "GSM"; Cell = "K70"; Value = 2801 }
    @{ Sheet = "GSM"; Cell = "M70"; Value = -2531 }
    @{ Sheet = "GSM"; Cell = "H73"; Value = 2801 }
    @{ Sheet = "GSM"; Cell = "I73"; Value = 2801 }
    @{ Sheet = "GSM"; Cell = "K73"; Value = 2801 }
    @{ Sheet = "GSM"; Cell = "M73"; Value = -1865 }
    @{ Sheet = "GSM"; Cell = "H80"; Value = 2930.889 }
    @{ Sheet = "GSM"; Cell = "I80"; Value = 2575.8 }
    @{ Sheet = "GSM"; Cell = "J80"; Value = 3374.75 }
    @{ Sheet = "GSM"; Cell = "K80"; Value = 2575.8 }
    @{ Sheet = "GSM"; Cell = "L80"; Value = 3374.75 }
    @{ Sheet = "GSM"; Cell = "M80"; Value = -1577.8 }
    @{ Sheet = "GSM"; Cell = "N80"; Value = -5370.75 }
    @{ Sheet = "GSM"; Cell = "H83"; Value = 2930.889 }
    @{ Sheet = "GSM"; Cell = "I83"; Value = 2575.8 }
    @{ Sheet = "GSM"; Cell = "J83"; Value = 3374.75 }
    @{ Sheet = "GSM"; Cell = "K83"; Value = 12879 }
    @{ Sheet = "GSM"; Cell = "L83"; Value = 16873.75 }
    @{ Sheet = "GSM"; Cell = "M83"; Value = -7887 }
    @{ Sheet = "GSM"; Cell = "N83"; Value = -26857.75 }
    @{ Sheet = "GSM"; Cell = "H94"; Value = 28776.637 }
    @{ Sheet = "GSM"; Cell = "J94"; Value = 29613.6 }
    @{ Sheet = "GSM"; Cell = "L94"; Value = 29613.6 }
    @{ Sheet = "GSM"; Cell = "N94"; Value = -30965.6 }
    @{ Sheet = "GSM"; Cell = "H107"; Value = 2200.2856 }
    @{ Sheet = "GSM"; Cell = "I107"; Value = 726.375 }
    @{ Sheet = "GSM"; Cell = "K107"; Value = 726.375 }
    @{ Sheet = "GSM"; Cell = "M107"; Value = 1193.625 }
    @{ Sheet = "LTW"; Cell = "H7"; Value = 7583.1055 }
    @{ Sheet = "LTW"; Cell = "I7"; Value = 7198.25 }
    @{ Sheet = "LTW"; Cell = "K7"; Value = 7198.25 }
    @{ Sheet = "LTW"; Cell = "M7"; Value = -7086.25 }
    @{ Sheet = "LTW"; Cell = "H46"; Value = 2322 }
    @{ Sheet = "LTW"; Cell = "I46"; Value = 992.25 }
    @{ Sheet = "LTW"; Cell = "K46"; Value = 992.25 }
    @{ Sheet = "LTW"; Cell = "M46"; Value = -804.25 }
    @{ Sheet = "LTW"; Cell = "H61"; Value = 1624.25 }
    @{ Sheet = "LTW"; Cell = "J61"; Value = 1999 }
    @{ Sheet = "LTW"; Cell = "L61"; Value = 1999 }
    @{ Sheet = "LTW"; Cell = "N61"; Value = -2403 }
    @{ Sheet = "LTW"; Cell = "H113"; Value = 1624.25 }
    @{ Sheet = "LTW"; Cell = "J113"; Value = 1999 }
    @{ Sheet = "LTW"; Cell = "L113"; Value = 1999 }
    @{ Sheet = "LTW"; Cell = "N113"; Value = -6339 }
    @{ Sheet = "LTW"; Cell = "H126"; Value = 7583.1055 }
    @{ Sheet = "LTW"; Cell = "I126"; Value = 7198.25 }
    @{ Sheet = "LTW"; Cell = "K126"; Value = 21594.75 }
    @{ Sheet = "LTW"; Cell = "M126"; Value = -19124.75 }
    @{ Sheet = "WVR"; Cell = "H132"; Value = 2129.1 }
    @{ Sheet = "WVR"; Cell = "I132"; Value = 2129.1 }
    @{ Sheet = "WVR"; Cell = "K132"; Value = 6387.299999999999 }
    @{ Sheet = "WVR"; Cell = "M132"; Value = -3857.299999999999 }
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    $ws.Range($edit.Cell).Value = $edit.Value
}

Write-Output "Applied $($edits.Count) cell updates."